# Add new antenna/IP rows to the scheme sheet, and renumber a couple of
# "Ant. N" labels whose sequence shifted because a row was inserted in the
# middle of them.
#
# Insertions are applied top-to-bottom using each new row's FINAL row
# number. Since every earlier insert has already been applied (and rows
# above the current insertion point are therefore already in their final
# positions), using the final row index for each subsequent Insert() call
# lands it in the correct spot without needing any further offset math.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows at the very top: "Company Ant. 1/2" / 169.254.1.1/.2 ---
$ws.Rows("1:2").Insert()
$ws.Range("A1").Value = "Company Ant. 1"
$ws.Range("B1").Value = "169.254.1.1"
$ws.Range("A2").Value = "Company Ant. 2"
$ws.Range("B2").Value = "169.254.1.2"

# --- New row: "V 2282 Ant. 8" / 169.254.1.27, ends up as row 14 ---
$ws.Rows(14).Insert()
$ws.Range("A14").Value = "V 2282 Ant. 8"
$ws.Range("B14").Value = "169.254.1.27"

# --- New row: "Mobinil Dina Farms Ant. 2" / 169.254.1.28, ends up as row 32 ---
$ws.Rows(32).Insert()
$ws.Range("A32").Value = "Mobinil Dina Farms Ant. 2"
$ws.Range("B32").Value = "169.254.1.28"

# The existing "Mobinil Dina Farms" row just below (IP 169.254.1.29) is
# renumbered from Ant. 2 to Ant. 3 now that a new Ant. 2 exists above it.
$ws.Range("A33").Value = "Mobinil Dina Farms Ant. 3"

# --- New row: "Mobinil Dina Farms Ant. 2" / 169.254.1.31, ends up as row 34 ---
$ws.Rows(34).Insert()
$ws.Range("A34").Value = "Mobinil Dina Farms Ant. 2"
$ws.Range("B34").Value = "169.254.1.31"

# --- New row: "4117 Ant. 2" / 169.254.1.32, ends up as row 38 ---
$ws.Rows(38).Insert()
$ws.Range("A38").Value = "4117 Ant. 2"
$ws.Range("B38").Value = "169.254.1.32"

# The existing "4117 Ant. 2/3/4" rows just below are bumped up to
# Ant. 3/4/5 now that a new Ant. 2 exists above them.
$ws.Range("A39").Value = "4117 Ant. 3"
$ws.Range("A40").Value = "4117 Ant. 4"
$ws.Range("A41").Value = "4117 Ant. 5"

# --- Match the view/selection state captured in the saved workbook ---
$ws.Range("B34").Select()
